$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($range, [string]$text)
    # Force the cell to be stored as text even when the string looks like a number,
    # then strip the number-format style back off so no stray style index is left behind.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Set-Cell {
    param($range, [string]$text)
    $range.Value = $text
}

# --- Row 8: was split "Documento sólo en Sofía" -> now a full matched row ---
Set-Cell      ($ws.Cells.Item(8, 5))  "CC"                              # E8
Set-TextCell  ($ws.Cells.Item(8, 7))  "167393272"                       # G8
Set-TextCell  ($ws.Cells.Item(8, 8))  "167393272"                       # H8
Set-Cell      ($ws.Cells.Item(8, 9))  "JOAQUIN BERNARDO LOPEZ PEREA"    # I8
Set-Cell      ($ws.Cells.Item(8, 11)) "VERDADERO"                       # K8

# --- Row 9: was split "Documento sólo en Instructores" -> now a full discrepancy row ---
Set-Cell      ($ws.Cells.Item(9, 1))  3031278                                      # A9
Set-Cell      ($ws.Cells.Item(9, 2))  "COMPLEMENTARIA"                             # B9
Set-Cell      ($ws.Cells.Item(9, 3))  "CURSO ESPECIAL"                             # C9
Set-Cell      ($ws.Cells.Item(9, 4))  "ELABORACION DE PRODUCTOS DE REPOSTERIA."    # D9
Set-Cell      ($ws.Cells.Item(9, 6))  "CC"                                         # F9
Set-TextCell  ($ws.Cells.Item(9, 7))  "18464762"                                   # G9
Set-TextCell  ($ws.Cells.Item(9, 8))  "18464762"                                   # H9
Set-Cell      ($ws.Cells.Item(9, 9))  "LUIS EDUARDO OLIVEROS RAMIRES"              # I9
Set-Cell      ($ws.Cells.Item(9, 10)) "LUIS EDUARDO OLIVEROS RAMIREZ"              # J9
Set-Cell      ($ws.Cells.Item(9, 11)) "FALSO - Discrepancia en Nombre: Instructores (LUIS EDUARDO OLIVEROS RAMIRES) vs Sofía (LUIS EDUARDO OLIVEROS RAMIREZ)" # K9

# --- Row 10 ---
Set-TextCell  ($ws.Cells.Item(10, 7))  "250188492"                      # G10
Set-TextCell  ($ws.Cells.Item(10, 8))  "250188492"                      # H10
Set-Cell      ($ws.Cells.Item(10, 9))  "NORA MARIA BLANDON SERNA"       # I10
Set-Cell      ($ws.Cells.Item(10, 10)) "NORA MARIA BLANDON SERNA"       # J10
Set-Cell      ($ws.Cells.Item(10, 11)) "VERDADERO"                      # K10

# --- Row 11 ---
Set-TextCell  ($ws.Cells.Item(11, 7))  "31642764"                       # G11
Set-TextCell  ($ws.Cells.Item(11, 8))  "31642764"                       # H11
Set-Cell      ($ws.Cells.Item(11, 9))  "INGRI YEANA ESCOBAR "           # I11
Set-Cell      ($ws.Cells.Item(11, 10)) "INGRI YEANA ESCOBAR "           # J11

# --- Row 12 ---
Set-TextCell  ($ws.Cells.Item(12, 7))  "31863992"                       # G12
Set-TextCell  ($ws.Cells.Item(12, 8))  "31863992"                       # H12
Set-Cell      ($ws.Cells.Item(12, 9))  "ALBA INES ZULETA JARAMILLO"     # I12
Set-Cell      ($ws.Cells.Item(12, 10)) "ALBA INES ZULETA JARAMILLO"     # J12

# --- Row 13 ---
Set-TextCell  ($ws.Cells.Item(13, 7))  "327311712"                            # G13
Set-TextCell  ($ws.Cells.Item(13, 8))  "327311712"                            # H13
Set-Cell      ($ws.Cells.Item(13, 9))  "ZAILAN DEL CARMEN CALDERON LOCARNO"   # I13
Set-Cell      ($ws.Cells.Item(13, 10)) "ZAILAN DEL CARMEN CALDERON LOCARNO"   # J13

# --- Row 14: was split "Documento sólo en Sofía" -> now a full matched row ---
Set-Cell      ($ws.Cells.Item(14, 5))  "CC"                              # E14
Set-TextCell  ($ws.Cells.Item(14, 7))  "36380685"                        # G14
Set-TextCell  ($ws.Cells.Item(14, 8))  "36380685"                        # H14
Set-Cell      ($ws.Cells.Item(14, 9))  "GLORIA YENNY CASTILLO ESPAÑA"    # I14
Set-Cell      ($ws.Cells.Item(14, 10)) "GLORIA YENNY CASTILLO ESPAÑA"    # J14
Set-Cell      ($ws.Cells.Item(14, 11)) "VERDADERO"                       # K14

# --- Row 15: was split "Documento sólo en Instructores" -> now a full matched row ---
Set-Cell      ($ws.Cells.Item(15, 1))  3031278                                      # A15
Set-Cell      ($ws.Cells.Item(15, 2))  "COMPLEMENTARIA"                             # B15
Set-Cell      ($ws.Cells.Item(15, 3))  "CURSO ESPECIAL"                             # C15
Set-Cell      ($ws.Cells.Item(15, 4))  "ELABORACION DE PRODUCTOS DE REPOSTERIA."    # D15
Set-Cell      ($ws.Cells.Item(15, 6))  "CC"                                         # F15
Set-TextCell  ($ws.Cells.Item(15, 7))  "38553002"                                   # G15
Set-TextCell  ($ws.Cells.Item(15, 8))  "38553002"                                   # H15
Set-Cell      ($ws.Cells.Item(15, 9))  "PAOLA ANDREA CASTILLO ALZATE"               # I15
Set-Cell      ($ws.Cells.Item(15, 10)) "PAOLA ANDREA CASTILLO ALZATE"               # J15
Set-Cell      ($ws.Cells.Item(15, 11)) "VERDADERO"                                  # K15

# --- Row 16 ---
Set-TextCell  ($ws.Cells.Item(16, 7))  "42057870"                       # G16
Set-TextCell  ($ws.Cells.Item(16, 8))  "42057870"                       # H16
Set-Cell      ($ws.Cells.Item(16, 9))  "MARTHA LUCIA ALZATE GOMEZ"      # I16
Set-Cell      ($ws.Cells.Item(16, 10)) "MARTHA LUCIA ALZATE GOMEZ"      # J16

# --- Row 17 ---
Set-TextCell  ($ws.Cells.Item(17, 7))  "66655884"                       # G17
Set-TextCell  ($ws.Cells.Item(17, 8))  "66655884"                       # H17
Set-Cell      ($ws.Cells.Item(17, 9))  "MARIA DEL CARMEN CERON BEDOYA"  # I17
Set-Cell      ($ws.Cells.Item(17, 10)) "MARIA DEL CARMEN CERON BEDOYA"  # J17

# --- Row 18 ---
Set-Cell      ($ws.Cells.Item(18, 6))  "TI"                             # F18
Set-TextCell  ($ws.Cells.Item(18, 7))  "66756576"                       # G18
Set-TextCell  ($ws.Cells.Item(18, 8))  "66756576"                       # H18
Set-Cell      ($ws.Cells.Item(18, 9))  "MARTHA ISABEL BENAVIDES ACOSTA" # I18
Set-Cell      ($ws.Cells.Item(18, 10)) "MARTHA ISABEL BENAVIDES ACOSTA" # J18
Set-Cell      ($ws.Cells.Item(18, 11)) "FALSO - Discrepancia en Tipo de Documento: Instructores (CC) vs Sofía (TI)" # K18

# --- Row 19 ---
Set-TextCell  ($ws.Cells.Item(19, 7))  "66767235"                       # G19
Set-TextCell  ($ws.Cells.Item(19, 8))  "66767235"                       # H19
Set-Cell      ($ws.Cells.Item(19, 9))  "CLAUDIA PATRICIA ARCE ESCOBAR"  # I19
Set-Cell      ($ws.Cells.Item(19, 10)) "CLAUDIA PATRICIA ARCE ESCOBAR"  # J19

# --- Rows 20 and 21 no longer exist: delete them entirely ---
$ws.Rows("20:21").Delete()
